$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F header "time_taken", copying the formatting of the
# existing header cell E1 (bold font, border, centered alignment)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Populate the time_taken values for data rows 2-8
$times = @(
    "2021-10-05 10:50:51.939006",
    "2021-10-05 10:50:51.939018",
    "2021-10-05 10:50:51.939022",
    "2021-10-05 10:50:51.939025",
    "2021-10-05 10:50:51.939028",
    "2021-10-05 10:50:51.939031",
    "2021-10-05 10:50:51.939035"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
